$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-30 Tuesday" "2024-07-31 Wednesday"

Replace-Text "180×5=900" "779×6=4674"
Replace-Text "138×2=276" "577×2=1154"
Replace-Text "857×9=7713" "310×5=1550"
Replace-Text "162×8=1296" "826×8=6608"
Replace-Text "520×8=4160" "266×5=1330"

Replace-Text "107×2=214" "914×7=6398"
Replace-Text "979×2=1958" "301×5=1505"
Replace-Text "188×2=376" "530×6=3180"
Replace-Text "488×4=1952" "133×6=798"
Replace-Text "546×5=2730" "574×5=2870"

Replace-Text "726×4=2904" "449×8=3592"
Replace-Text "698×2=1396" "878×9=7902"
Replace-Text "311×8=2488" "157×9=1413"
Replace-Text "451×4=1804" "664×4=2656"
Replace-Text "468×7=3276" "453×8=3624"

Replace-Text "541×5=2705" "242×6=1452"
Replace-Text "807×6=4842" "828×7=5796"
Replace-Text "252×6=1512" "583×8=4664"
Replace-Text "397×7=2779" "364×5=1820"
Replace-Text "828×6=4968" "843×9=7587"

Replace-Text "774×9=6966" "650×8=5200"
Replace-Text "821×6=4926" "567×8=4536"
Replace-Text "967×6=5802" "655×8=5240"
Replace-Text "993×6=5958" "507×7=3549"
Replace-Text "792×5=3960" "650×6=3900"
